# Damped hind leg calculation - start formatting measured new leg property data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Measured values per limb (entered first)
$ws.Range("H8").Value = "limb 1"
$ws.Range("I8").Value = 190.04
$ws.Range("J8").Value = 130.77000000000001
$ws.Range("K8").Value = 192.72

$ws.Range("H9").Value = "limb 2"
$ws.Range("I9").Value = 267.99700000000001
$ws.Range("J9").Value = 171.35
$ws.Range("K9").Value = 166.26

$ws.Range("H10").Value = "limb 3"
$ws.Range("I10").Value = 162.02000000000001
$ws.Range("J10").Value = 53.49
$ws.Range("K10").Value = 35.33

# Column labels for the three limbs (entered after the data)
$ws.Range("I7").Value = "L"
$ws.Range("J7").Value = "R"
$ws.Range("K7").Value = "m"

# Header for the new measured data block, merged across H6:K6, centered (entered last)
$ws.Range("H6").Value = "measured new leg properties"
$ws.Range("H6:K6").HorizontalAlignment = -4108  # xlCenter
$ws.Range("H6:K6").Merge() | Out-Null

# Leave the active selection where the author left off
$ws.Range("K11").Select() | Out-Null
